$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cells (Wins / Losses / Ties) in AD1:AF1, reusing the
# same header style (bold font, border, centered) already applied to the
# rest of row 1 by copying the format from an existing header cell.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins=73, Losses=89, Ties=0) for every player
# row (rows 2 through 59).
$lastRow = 59
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 73
    $ws.Cells.Item($r, 31).Value = 89
    $ws.Cells.Item($r, 32).Value = 0
}
